# Update the dSF column (F) values per repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 7
    4  = -1
    5  = -6
    6  = -6
    7  = 5
    8  = 2
    9  = 2
    10 = -3
    13 = -4
    14 = 2
    15 = -1
    16 = -1
    17 = -3
    18 = -2
    19 = 1
    20 = 7
    21 = 2
    22 = 2
    23 = -3
    25 = 2
    26 = -5
    27 = 0
    28 = 1
    29 = -4
    30 = -8
    31 = 1
    32 = 2
    33 = -5
    34 = 7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
